$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds the "changed" date serial 45203 (2023-10-04)
# for every data row (rows 2-340). Update it to 45204 (2023-10-05).
for ($row = 2; $row -le 340; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45203) {
        $cell.Value = 45204
    }
}
